# Update "paises" COVID dashboard: countries & provincias Spain refresh.
# Reflects a later data pull (27 Mar 2020, 23:53) -- updated totals for a
# handful of countries, including a few whose running total overtook a
# neighboring country in the ranking (so that neighbor's row shifts down
# by one while keeping its own numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    $col = 1
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col++
    }
}

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 23:53"

# Straightforward number refreshes (no ranking change)
Set-Row 4   @("Estados Unidos", 101652, 16217, 2465, 97599, 2463, 293, 1588)
Set-Row 18  @("Canada", 4689, 646, 258, 4378, 120, 14, 53)
Set-Row 34  @("Pakistan", 1373, 172, 23, 1339, 7, 2, 11)
Set-Row 49  @("Peru", 635, 55, 16, 608, 21, 2, 11)
Set-Row 138 @("Togo", 25, 1, 1, 23, 0, 1, 1)

# Bosnia y Herzegovina overtakes Jordania -> rows 79/80 swap places;
# Bosnia gets fresh numbers, Jordania keeps its previous ones.
Set-Row 79 @("Bosnia y Herzegovina", 237, 46, 5, 228, 1, 1, 4)
Set-Row 80 @("Jordania", 235, 23, 18, 216, 0, 1, 1)

# Aruba overtakes Guam, Kenia, Polinesia Francesa and Isla de Man ->
# rows 129-133 shift down one place; Aruba gets fresh numbers, the rest
# keep their previous ones, just one row lower.
Set-Row 129 @("Aruba", 33, 5, 1, 32, 0, 0, 0)
Set-Row 130 @("Guam", 32, 0, 0, 31, 0, 0, 1)
Set-Row 131 @("Kenia", 31, 0, 1, 29, 2, 0, 1)
Set-Row 132 @("Polinesia Francesa", 30, 0, 0, 30, 0, 0, 0)
Set-Row 133 @("Isla de Man", 29, 3, 0, 29, 0, 0, 0)

# Bahamas overtakes Niger, Groenlandia and Suazilandia -> rows 155-158
# shift down one place; Bahamas gets fresh numbers, the rest keep their
# previous ones, just one row lower.
Set-Row 155 @("Bahamas", 10, 1, 1, 9, 0, 0, 0)
Set-Row 156 @("Niger", 10, 0, 0, 9, 0, 0, 1)
Set-Row 157 @("Groenlandia", 10, 4, 2, 8, 0, 0, 0)
Set-Row 158 @("Suazilandia", 9, 3, 0, 9, 0, 0, 0)
